{"js": "// 1) Insert a new \"ListBullet\" paragraph with the new docente right after the\n//    \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph.\nconst body = context.document.body;\nconst heading = body.search(\"Docente(s) Respons\u00e1vel(eis) \", { matchCase: true });\nheading.load(\"items\");\nawait context.sync();\n\nif (heading.items.length > 0) {\n  const headingPara = heading.items[0].paragraphs.getFirst();\n  const newPara = headingPara.insertParagraph(\"5817181 - Valdeir Arantes\", \"After\");\n  newPara.style = \"ListBullet\";\n}\n\n// 2) Update the LOT2058 requirement line's course-name text, keeping the\n//    trailing \"(Requisito fraco)\" and the line break that follows it intact.\nconst reqHits = body.search(\"LOT2058 -  Engenharia Econ\u00f4mica  (Requisito fraco)\", { matchCase: true });\nreqHits.load(\"items\");\nawait context.sync();\n\nif (reqHits.items.length > 0) {\n  reqHits.items[0].insertText(\"LOT2058 -  An\u00e1lise T\u00e9cnico (Requisito fraco)\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new \"ListBullet\" paragraph with the new docente right after the\n#    \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Docente(s) Respons\u00e1vel(eis) \") {\n        $p.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newPara.Range.Text = \"5817181 - Valdeir Arantes\"\n        $newPara.Style = \"ListBullet\"\n        break\n    }\n}\n\n# 2) Update the LOT2058 requirement line's course-name text, keeping the\n#    trailing \"(Requisito fraco)\" and the line break that follows it intact.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"LOT2058 -  Engenharia Econ\u00f4mica  (Requisito fraco)\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"LOT2058 -  An\u00e1lise T\u00e9cnico (Requisito fraco)\",\n    2\n)\n"}
